# Add detailed data for Data Collection rows 92-106 (new bone specimens, Day 360 series)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Collection")

# Row 92
$ws.Range("B92").Value = "FEMUR"
$ws.Range("C92").Value = 1
$ws.Range("D92").Value = 1
$ws.Range("E92").Value = "OBLIQUE"
$ws.Range("F92").Value = "COMMINUTED"
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = "INTERMEDIATE"
$ws.Range("I92").Value = "CURVED/JAGGED"
$ws.Range("J92").Value = "N/A"
$ws.Range("K92").Value = "N/A"

# Row 93
$ws.Range("B93").Value = "FEMUR"
$ws.Range("C93").Value = 1
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = "INTERMEDIATE"
$ws.Range("F93").Value = "SEGMENTAL"
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = "SMOOTH"
$ws.Range("I93").Value = "CURVED/JAGGED"
$ws.Range("J93").Value = "N/A"
$ws.Range("K93").Value = "N/A"

# Row 94
$ws.Range("B94").Value = "HUMERUS"
$ws.Range("C94").Value = 1
$ws.Range("D94").Value = 1
$ws.Range("E94").Value = "OBLIQUE"
$ws.Range("F94").Value = "COMMINUTED"
$ws.Range("G94").Value = 10
$ws.Range("H94").Value = "SMOOTH"
$ws.Range("I94").Value = "CURVED/JAGGED"
$ws.Range("J94").Value = "N/A"
$ws.Range("K94").Value = "N/A"

# Row 95
$ws.Range("B95").Value = "FEMUR"
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = "OBLIQUE"
$ws.Range("F95").Value = "COMMINUTED"
$ws.Range("G95").Value = 4
$ws.Range("H95").Value = "INTERMEDIATE"
$ws.Range("I95").Value = "CURVED/JAGGED"
$ws.Range("J95").Value = "BONE HACKLE "
$ws.Range("K95").Value = "<25%"

# Row 96
$ws.Range("B96").Value = "FEMUR"
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 1
$ws.Range("E96").Value = "OBLIQUE"
$ws.Range("F96").Value = "COMMINUTED"
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = "INTERMEDIATE"
$ws.Range("I96").Value = "CURVED/JAGGED"
$ws.Range("J96").Value = "N/A"
$ws.Range("K96").Value = "N/A"

# Row 97
$ws.Range("B97").Value = "TIBIA"
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1
$ws.Range("E97").Value = "OBLIQUE"
$ws.Range("F97").Value = "COMMINUTED"
$ws.Range("G97").Value = 6
$ws.Range("H97").Value = "SMOOTH"
$ws.Range("I97").Value = "CURVED/JAGGED"
$ws.Range("J97").Value = "N/A"
$ws.Range("K97").Value = "N/A"

# Row 98
$ws.Range("B98").Value = "TIBIA"
$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 1
$ws.Range("E98").Value = "INTERMEDIATE"
$ws.Range("F98").Value = "BUTTERFLY"
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = "INTERMEDIATE"
$ws.Range("I98").Value = "CURVED/JAGGED"
$ws.Range("J98").Value = "N/A"
$ws.Range("K98").Value = "N/A"

# Row 99
$ws.Range("B99").Value = "TIBIA"
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 1
$ws.Range("E99").Value = "OBLIQUE"
$ws.Range("F99").Value = "COMMINUTED"
$ws.Range("G99").Value = 5
$ws.Range("H99").Value = "INTERMEDIATE"
$ws.Range("I99").Value = "CURVED/JAGGED"
$ws.Range("J99").Value = "COMPRESSION CURL"
$ws.Range("K99").Value = "<25%"

# Row 100
$ws.Range("B100").Value = "FEMUR"
$ws.Range("C100").Value = 1
$ws.Range("D100").Value = 1
$ws.Range("E100").Value = "OBLIQUE"
$ws.Range("F100").Value = "COMMINUTED"
$ws.Range("G100").Value = 8
$ws.Range("H100").Value = "INTERMEDIATE"
$ws.Range("I100").Value = "CURVED/JAGGED"
$ws.Range("J100").Value = "N/A"
$ws.Range("K100").Value = "N/A"

# Row 101
$ws.Range("B101").Value = "RADIUS"
$ws.Range("C101").Value = 1
$ws.Range("D101").Value = 1
$ws.Range("E101").Value = "OBLIQUE"
$ws.Range("F101").Value = "COMMINUTED"
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = "INTERMEDIATE"
$ws.Range("I101").Value = "CURVED/JAGGED"
$ws.Range("J101").Value = "N/A"
$ws.Range("K101").Value = "N/A"

# Row 102
$ws.Range("B102").Value = "FEMUR"
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 1
$ws.Range("E102").Value = "INTERMEDIATE"
$ws.Range("F102").Value = "COMMINUTED"
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = "INTERMEDIATE"
$ws.Range("I102").Value = "CURVED/JAGGED"
$ws.Range("J102").Value = "N/A"
$ws.Range("K102").Value = "N/A"

# Row 103
$ws.Range("B103").Value = "TIBIA"
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 1
$ws.Range("E103").Value = "OBLIQUE"
$ws.Range("F103").Value = "COMMINUTED"
$ws.Range("G103").Value = 3
$ws.Range("H103").Value = "SMOOTH"
$ws.Range("I103").Value = "CURVED/JAGGED"
$ws.Range("J103").Value = "WAKE HACKLE"
$ws.Range("K103").Value = "N/A"

# Row 104
$ws.Range("B104").Value = "HUMERUS"
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = "N/A"
$ws.Range("F104").Value = "N/A"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = "N/A"
$ws.Range("I104").Value = "N/A"
$ws.Range("J104").Value = "N/A"
$ws.Range("K104").Value = "N/A"

# Row 105
$ws.Range("B105").Value = "HUMERUS"
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 1
$ws.Range("E105").Value = "N/A"
$ws.Range("F105").Value = "N/A"
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = "N/A"
$ws.Range("I105").Value = "N/A"
$ws.Range("J105").Value = "BONE HACKLE, WAKE HACKLE, ARREST RIDGE"
$ws.Range("K105").Value = "<25%, <25%, <25%"

# Row 106
$ws.Range("B106").Value = "HUMERUS"
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 1
$ws.Range("E106").Value = "INTERMEDIATE"
$ws.Range("F106").Value = "COMMINUTED"
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = "ROUGH"
$ws.Range("I106").Value = "CURVED/JAGGED"
$ws.Range("J106").Value = "N/A"
$ws.Range("K106").Value = "N/A"

# Re-apply the sheet-wide centered alignment style used throughout this table to the newly
# populated cells so they match the existing formatting.
$ws.Range("B92:K106").HorizontalAlignment = -4108

# Restore the view position / selection recorded for this sheet.
$ws.Activate()
[void]$ws.Range("L106").Select()
$excel.ActiveWindow.ScrollRow = 83
